$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 32.736679
$ws.Cells.Item(2, 8).Value = 98.210037
$ws.Cells.Item(2, 9).Value = 0.8346853755332739
$ws.Cells.Item(2, 10).Value = 0.834685375533274
$ws.Cells.Item(2, 13).Value = 90.25004833333332
$ws.Cells.Item(2, 14).Value = 270.750145
$ws.Cells.Item(2, 15).Value = 0.8928575650827933
$ws.Cells.Item(2, 16).Value = 0.8928575650827932
$ws.Cells.Item(2, 17).Value = 2954.486862022818
$ws.Cells.Item(2, 18).Value = 26590.38175820536
$ws.Cells.Item(2, 19).Value = 0.7452551520088559
$ws.Cells.Item(2, 20).Value = 0.7452551520088559
# Row 3
$ws.Cells.Item(3, 7).Value = 32.736679
$ws.Cells.Item(3, 8).Value = 98.210037
$ws.Cells.Item(3, 9).Value = 0.8346853755332739
$ws.Cells.Item(3, 10).Value = 0.834685375533274
$ws.Cells.Item(3, 15).Value = 0.05133510428912089
$ws.Cells.Item(3, 16).Value = 0.05133510428912089
$ws.Cells.Item(3, 17).Value = 169.8690777948607
$ws.Cells.Item(3, 18).Value = 1528.821700153746
$ws.Cells.Item(3, 19).Value = 0.04284866080160465
$ws.Cells.Item(3, 20).Value = 0.04284866080160465
# Row 4
$ws.Cells.Item(4, 7).Value = 32.736679
$ws.Cells.Item(4, 8).Value = 98.210037
$ws.Cells.Item(4, 9).Value = 0.8346853755332739
$ws.Cells.Item(4, 10).Value = 0.834685375533274
$ws.Cells.Item(4, 13).Value = 5.380476000000001
$ws.Cells.Item(4, 14).Value = 16.141428
$ws.Cells.Item(4, 15).Value = 0.05322987398968605
$ws.Cells.Item(4, 16).Value = 0.05322987398968604
$ws.Cells.Item(4, 17).Value = 176.138915679204
$ws.Cells.Item(4, 18).Value = 1585.250241112836
$ws.Cells.Item(4, 19).Value = 0.04443019736066995
$ws.Cells.Item(4, 20).Value = 0.04443019736066995
# Row 5
$ws.Cells.Item(5, 7).Value = 32.736679
$ws.Cells.Item(5, 8).Value = 98.210037
$ws.Cells.Item(5, 9).Value = 0.8346853755332739
$ws.Cells.Item(5, 10).Value = 0.834685375533274
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.2605293333333333
$ws.Cells.Item(5, 14).Value = 0.7815879999999999
$ws.Cells.Item(5, 15).Value = 0.002577456638399696
$ws.Cells.Item(5, 16).Value = 0.002577456638399696
$ws.Cells.Item(5, 17).Value = 8.528865155417334
$ws.Cells.Item(5, 18).Value = 76.75978639875599
$ws.Cells.Item(5, 19).Value = 0.002151365362143381
$ws.Cells.Item(5, 20).Value = 0.002151365362143381
# Row 6
$ws.Cells.Item(6, 9).Value = 0.001766029048926899
$ws.Cells.Item(6, 10).Value = 0.0017660290489269
$ws.Cells.Item(6, 13).Value = 90.25004833333332
$ws.Cells.Item(6, 14).Value = 270.750145
$ws.Cells.Item(6, 15).Value = 0.8928575650827933
$ws.Cells.Item(6, 16).Value = 0.8928575650827932
$ws.Cells.Item(6, 17).Value = 6.251109431109444
$ws.Cells.Item(6, 18).Value = 56.259984879985
$ws.Cells.Item(6, 19).Value = 0.001576812396490352
$ws.Cells.Item(6, 20).Value = 0.001576812396490353
# Row 7
$ws.Cells.Item(7, 9).Value = 0.001766029048926899
$ws.Cells.Item(7, 10).Value = 0.0017660290489269
$ws.Cells.Item(7, 15).Value = 0.05133510428912089
$ws.Cells.Item(7, 16).Value = 0.05133510428912089
$ws.Cells.Item(7, 19).Value = 0.00009065928540427934
$ws.Cells.Item(7, 20).Value = 0.00009065928540427937
# Row 8
$ws.Cells.Item(8, 9).Value = 0.001766029048926899
$ws.Cells.Item(8, 10).Value = 0.0017660290489269
$ws.Cells.Item(8, 13).Value = 5.380476000000001
$ws.Cells.Item(8, 14).Value = 16.141428
$ws.Cells.Item(8, 15).Value = 0.05322987398968605
$ws.Cells.Item(8, 16).Value = 0.05322987398968604
$ws.Cells.Item(8, 17).Value = 0.372675083156
$ws.Cells.Item(8, 18).Value = 3.354075748404
$ws.Cells.Item(8, 19).Value = 0.00009400550373650393
$ws.Cells.Item(8, 20).Value = 0.00009400550373650394
# Row 9
$ws.Cells.Item(9, 9).Value = 0.001766029048926899
$ws.Cells.Item(9, 10).Value = 0.0017660290489269
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.2605293333333333
$ws.Cells.Item(9, 14).Value = 0.7815879999999999
$ws.Cells.Item(9, 15).Value = 0.002577456638399696
$ws.Cells.Item(9, 16).Value = 0.002577456638399696
$ws.Cells.Item(9, 17).Value = 0.01804539058711111
$ws.Cells.Item(9, 18).Value = 0.162408515284
$ws.Cells.Item(9, 19).Value = 0.000004551863295763338
$ws.Cells.Item(9, 20).Value = 0.000004551863295763339
# Row 10
$ws.Cells.Item(10, 7).Value = 0.5119106666666666
$ws.Cells.Item(10, 8).Value = 1.535732
$ws.Cells.Item(10, 9).Value = 0.01305215923234471
$ws.Cells.Item(10, 10).Value = 0.01305215923234471
$ws.Cells.Item(10, 13).Value = 90.25004833333332
$ws.Cells.Item(10, 14).Value = 270.750145
$ws.Cells.Item(10, 15).Value = 0.8928575650827933
$ws.Cells.Item(10, 16).Value = 0.8928575650827932
$ws.Cells.Item(10, 17).Value = 46.19996240901555
$ws.Cells.Item(10, 18).Value = 415.7996616811399
$ws.Cells.Item(10, 19).Value = 0.0116537191112642
$ws.Cells.Item(10, 20).Value = 0.0116537191112642
# Row 11
$ws.Cells.Item(11, 7).Value = 0.5119106666666666
$ws.Cells.Item(11, 8).Value = 1.535732
$ws.Cells.Item(11, 9).Value = 0.01305215923234471
$ws.Cells.Item(11, 10).Value = 0.01305215923234471
$ws.Cells.Item(11, 15).Value = 0.05133510428912089
$ws.Cells.Item(11, 16).Value = 0.05133510428912089
$ws.Cells.Item(11, 17).Value = 2.656280218895111
$ws.Cells.Item(11, 18).Value = 23.906521970056
$ws.Cells.Item(11, 19).Value = 0.0006700339553906278
$ws.Cells.Item(11, 20).Value = 0.0006700339553906278
# Row 12
$ws.Cells.Item(12, 7).Value = 0.5119106666666666
$ws.Cells.Item(12, 8).Value = 1.535732
$ws.Cells.Item(12, 9).Value = 0.01305215923234471
$ws.Cells.Item(12, 10).Value = 0.01305215923234471
$ws.Cells.Item(12, 13).Value = 5.380476000000001
$ws.Cells.Item(12, 14).Value = 16.141428
$ws.Cells.Item(12, 15).Value = 0.05322987398968605
$ws.Cells.Item(12, 16).Value = 0.05322987398968604
$ws.Cells.Item(12, 17).Value = 2.754323056144
$ws.Cells.Item(12, 18).Value = 24.788907505296
$ws.Cells.Item(12, 19).Value = 0.0006947647912310264
$ws.Cells.Item(12, 20).Value = 0.0006947647912310264
# Row 13
$ws.Cells.Item(13, 7).Value = 0.5119106666666666
$ws.Cells.Item(13, 8).Value = 1.535732
$ws.Cells.Item(13, 9).Value = 0.01305215923234471
$ws.Cells.Item(13, 10).Value = 0.01305215923234471
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.2605293333333333
$ws.Cells.Item(13, 14).Value = 0.7815879999999999
$ws.Cells.Item(13, 15).Value = 0.002577456638399696
$ws.Cells.Item(13, 16).Value = 0.002577456638399696
$ws.Cells.Item(13, 17).Value = 0.1333677447128889
$ws.Cells.Item(13, 18).Value = 1.200309702416
$ws.Cells.Item(13, 19).Value = 0.00003364137445885676
$ws.Cells.Item(13, 20).Value = 0.00003364137445885676
# Row 14
$ws.Cells.Item(14, 7).Value = 5.902527666666667
$ws.Cells.Item(14, 8).Value = 17.707583
$ws.Cells.Item(14, 9).Value = 0.1504964361854544
$ws.Cells.Item(14, 10).Value = 0.1504964361854544
$ws.Cells.Item(14, 13).Value = 90.25004833333332
$ws.Cells.Item(14, 14).Value = 270.750145
$ws.Cells.Item(14, 15).Value = 0.8928575650827933
$ws.Cells.Item(14, 16).Value = 0.8928575650827932
$ws.Cells.Item(14, 17).Value = 532.7034072055038
$ws.Cells.Item(14, 18).Value = 4794.330664849535
$ws.Cells.Item(14, 19).Value = 0.1343718815661828
$ws.Cells.Item(14, 20).Value = 0.1343718815661828
# Row 15
$ws.Cells.Item(15, 7).Value = 5.902527666666667
$ws.Cells.Item(15, 8).Value = 17.707583
$ws.Cells.Item(15, 9).Value = 0.1504964361854544
$ws.Cells.Item(15, 10).Value = 0.1504964361854544
$ws.Cells.Item(15, 15).Value = 0.05133510428912089
$ws.Cells.Item(15, 16).Value = 0.05133510428912089
$ws.Cells.Item(15, 17).Value = 30.62793667602378
$ws.Cells.Item(15, 18).Value = 275.651430084214
$ws.Cells.Item(15, 19).Value = 0.007725750246721328
$ws.Cells.Item(15, 20).Value = 0.00772575024672133
# Row 16
$ws.Cells.Item(16, 7).Value = 5.902527666666667
$ws.Cells.Item(16, 8).Value = 17.707583
$ws.Cells.Item(16, 9).Value = 0.1504964361854544
$ws.Cells.Item(16, 10).Value = 0.1504964361854544
$ws.Cells.Item(16, 13).Value = 5.380476000000001
$ws.Cells.Item(16, 14).Value = 16.141428
$ws.Cells.Item(16, 15).Value = 0.05322987398968605
$ws.Cells.Item(16, 16).Value = 0.05322987398968604
$ws.Cells.Item(16, 17).Value = 31.758408449836
$ws.Cells.Item(16, 18).Value = 285.825676048524
$ws.Cells.Item(16, 19).Value = 0.008010906334048565
$ws.Cells.Item(16, 20).Value = 0.008010906334048567
# Row 17
$ws.Cells.Item(17, 7).Value = 5.902527666666667
$ws.Cells.Item(17, 8).Value = 17.707583
$ws.Cells.Item(17, 9).Value = 0.1504964361854544
$ws.Cells.Item(17, 10).Value = 0.1504964361854544
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.2605293333333333
$ws.Cells.Item(17, 14).Value = 0.7815879999999999
$ws.Cells.Item(17, 15).Value = 0.002577456638399696
$ws.Cells.Item(17, 16).Value = 0.002577456638399696
$ws.Cells.Item(17, 17).Value = 1.537781597978222
$ws.Cells.Item(17, 18).Value = 13.840034381804
$ws.Cells.Item(17, 19).Value = 0.0003878980385016958
$ws.Cells.Item(17, 20).Value = 0.0003878980385016958

$wb.Save()
